$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 25888
$ws.Range("E2").Value = 1403
$ws.Range("F2").Value = 1403
$ws.Range("G2").Value = 1139
$ws.Range("H2").Value = 900
$ws.Range("I2").Value = 928
$ws.Range("J2").Value = -28
$ws.Range("K2").Value = 18148
$ws.Range("L2").Value = 10303
$ws.Range("M2").Value = 7845
$ws.Range("N2").Value = 7579
$ws.Range("O2").Value = 266
$ws.Range("P2").Value = 360
$ws.Range("Q2").Value = 1627
$ws.Range("R2").Value = -1613
$ws.Range("S2").Value = -463
$ws.Range("T2").Value = 1075
$ws.Range("U2").Value = 553
$ws.Range("V2").Value = 6229
$ws.Range("W2").Value = 5.42
$ws.Range("X2").Value = 3.48
$ws.Range("Y2").Value = 12.87
$ws.Range("Z2").Value = 5.07
$ws.Range("AA2").Value = 131.34
$ws.Range("AB2").Value = 2122.78
$ws.Range("AC2").Value = 2576
$ws.Range("AD2").Value = 12.97
$ws.Range("AE2").Value = 21042
$ws.Range("AF2").Value = 1.59
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 0.9
$ws.Range("AI2").Value = 11.66
$ws.Range("AJ2").Value = 34411575
$ws.Range("D3").Value = 26350
$ws.Range("E3").Value = 1099
$ws.Range("F3").Value = 1099
$ws.Range("G3").Value = 777
$ws.Range("H3").Value = 472
$ws.Range("I3").Value = 583
$ws.Range("J3").Value = -111
$ws.Range("K3").Value = 21046
$ws.Range("L3").Value = 12861
$ws.Range("M3").Value = 8185
$ws.Range("N3").Value = 7934
$ws.Range("O3").Value = 252
$ws.Range("P3").Value = 360
$ws.Range("Q3").Value = 924
$ws.Range("R3").Value = -1875
$ws.Range("S3").Value = 1946
$ws.Range("T3").Value = 1849
$ws.Range("U3").Value = -925
$ws.Range("V3").Value = 8242
$ws.Range("W3").Value = 4.17
$ws.Range("X3").Value = 1.79
$ws.Range("Y3").Value = 7.51
$ws.Range("Z3").Value = 2.41
$ws.Range("AA3").Value = 157.12
$ws.Range("AB3").Value = 2223.39
$ws.Range("AC3").Value = 1618
$ws.Range("AD3").Value = 20.15
$ws.Range("AE3").Value = 22027
$ws.Range("AF3").Value = 1.48
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 1.23
$ws.Range("AI3").Value = 24.74
$ws.Range("AJ3").Value = 34411575
$ws.Range("D4").Value = 28550
$ws.Range("E4").Value = 1111
$ws.Range("F4").Value = 1111
$ws.Range("G4").Value = 819
$ws.Range("H4").Value = 601
$ws.Range("I4").Value = 646
$ws.Range("J4").Value = -46
$ws.Range("K4").Value = 23085
$ws.Range("L4").Value = 14412
$ws.Range("M4").Value = 8673
$ws.Range("N4").Value = 8482
$ws.Range("O4").Value = 191
$ws.Range("P4").Value = 360
$ws.Range("Q4").Value = 548
$ws.Range("R4").Value = -2203
$ws.Range("S4").Value = 967
$ws.Range("T4").Value = 1172
$ws.Range("U4").Value = -625
$ws.Range("V4").Value = 9576
$ws.Range("W4").Value = 3.89
$ws.Range("X4").Value = 2.1
$ws.Range("Y4").Value = 7.87
$ws.Range("Z4").Value = 2.72
$ws.Range("AA4").Value = 166.18
$ws.Range("AB4").Value = 2353.67
$ws.Range("AC4").Value = 1794
$ws.Range("AD4").Value = 14.72
$ws.Range("AE4").Value = 23549
$ws.Range("AF4").Value = 1.12
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 1.52
$ws.Range("AI4").Value = 22.32
$ws.Range("AJ4").Value = 34411575
$ws.Range("D5").Value = 29688
$ws.Range("E5").Value = 967
$ws.Range("F5").Value = 967
$ws.Range("G5").Value = 835
$ws.Range("H5").Value = 518
$ws.Range("I5").Value = 527
$ws.Range("J5").Value = -10
$ws.Range("K5").Value = 22524
$ws.Range("L5").Value = 13680
$ws.Range("M5").Value = 8844
$ws.Range("N5").Value = 8710
$ws.Range("O5").Value = 135
$ws.Range("P5").Value = 360
$ws.Range("Q5").Value = 1662
$ws.Range("R5").Value = -677
$ws.Range("S5").Value = -1094
$ws.Range("T5").Value = 1390
$ws.Range("U5").Value = 272
$ws.Range("V5").Value = 8511
$ws.Range("W5").Value = 3.26
$ws.Range("X5").Value = 1.74
$ws.Range("Y5").Value = 6.13
$ws.Range("Z5").Value = 2.27
$ws.Range("AA5").Value = 154.67
$ws.Range("AB5").Value = 2475.36
$ws.Range("AC5").Value = 1464
$ws.Range("AD5").Value = 18.82
$ws.Range("AE5").Value = 24181
$ws.Range("AF5").Value = 1.14
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 1.81
$ws.Range("AI5").Value = 34.18
$ws.Range("AJ5").Value = 34648025
$ws.Range("D6").Value = 29568
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = 1202
$ws.Range("G6").Value = 1191
$ws.Range("H6").Value = 654
$ws.Range("I6").Value = 649
$ws.Range("K6").Value = 21990
$ws.Range("L6").Value = 12703
$ws.Range("M6").Value = 9287
$ws.Range("N6").Value = 9151
$ws.Range("P6").Value = 360
$ws.Range("Q6").Value = 1216
$ws.Range("R6").Value = -234
$ws.Range("S6").Value = -752
$ws.Range("T6").Value = 1456
$ws.Range("U6").Value = -240
$ws.Range("V6").Value = 7952
$ws.Range("W6").Value = 4.06
$ws.Range("X6").Value = 2.21
$ws.Range("Y6").Value = 7.27
$ws.Range("Z6").Value = 2.94
$ws.Range("AA6").Value = 136.78
$ws.Range("AB6").Value = 2626.23
$ws.Range("AC6").Value = 1801
$ws.Range("AD6").Value = 13.91
$ws.Range("AE6").Value = 25405
$ws.Range("AF6").Value = 0.99
$ws.Range("AG6").Value = 550
$ws.Range("AH6").Value = 2.2
$ws.Range("AI6").Value = 30.56
$ws.Range("AJ6").Value = 34648025
$ws.Range("D7").Value = 29740
$ws.Range("E7").Value = 1312
$ws.Range("G7").Value = 1270
$ws.Range("H7").Value = 931
$ws.Range("I7").Value = 927
$ws.Range("K7").Value = 22414
$ws.Range("L7").Value = 12412
$ws.Range("M7").Value = 10002
$ws.Range("N7").Value = 9861
$ws.Range("P7").Value = 360
$ws.Range("Q7").Value = 1804
$ws.Range("R7").Value = -968
$ws.Range("S7").Value = -568
$ws.Range("T7").Value = 1122
$ws.Range("U7").Value = 649
$ws.Range("W7").Value = 4.41
$ws.Range("X7").Value = 3.13
$ws.Range("Y7").Value = 9.75
$ws.Range("Z7").Value = 4.2
$ws.Range("AA7").Value = 124.1
$ws.Range("AC7").Value = 2574
$ws.Range("AD7").Value = 8.140000000000001
$ws.Range("AE7").Value = 27377
$ws.Range("AF7").Value = 0.77
$ws.Range("AG7").Value = 550
$ws.Range("AH7").Value = 2.63
$ws.Range("AI7").Value = 20.56
$ws.Range("D8").Value = 30749
$ws.Range("E8").Value = 1437
$ws.Range("G8").Value = 1301
$ws.Range("H8").Value = 953
$ws.Range("I8").Value = 950
$ws.Range("K8").Value = 23133
$ws.Range("L8").Value = 12403
$ws.Range("M8").Value = 10730
$ws.Range("N8").Value = 10583
$ws.Range("P8").Value = 360
$ws.Range("Q8").Value = 1832
$ws.Range("R8").Value = -1021
$ws.Range("S8").Value = -419
$ws.Range("T8").Value = 1118
$ws.Range("U8").Value = 671
$ws.Range("W8").Value = 4.67
$ws.Range("X8").Value = 3.1
$ws.Range("Y8").Value = 9.289999999999999
$ws.Range("Z8").Value = 4.18
$ws.Range("AA8").Value = 115.59
$ws.Range("AC8").Value = 2638
$ws.Range("AD8").Value = 7.94
$ws.Range("AE8").Value = 29383
$ws.Range("AF8").Value = 0.71
$ws.Range("AG8").Value = 562
$ws.Range("AH8").Value = 2.68
$ws.Range("AI8").Value = 20.52
$ws.Range("D9").Value = 31924
$ws.Range("E9").Value = 1581
$ws.Range("G9").Value = 1454
$ws.Range("H9").Value = 1065
$ws.Range("I9").Value = 1060
$ws.Range("K9").Value = 23967
$ws.Range("L9").Value = 12409
$ws.Range("M9").Value = 11558
$ws.Range("N9").Value = 11405
$ws.Range("P9").Value = 360
$ws.Range("Q9").Value = 1963
$ws.Range("R9").Value = -1062
$ws.Range("S9").Value = -415
$ws.Range("T9").Value = 1126
$ws.Range("U9").Value = 753
$ws.Range("W9").Value = 4.95
$ws.Range("X9").Value = 3.34
$ws.Range("Y9").Value = 9.640000000000001
$ws.Range("Z9").Value = 4.52
$ws.Range("AA9").Value = 107.36
$ws.Range("AC9").Value = 2944
$ws.Range("AD9").Value = 7.12
$ws.Range("AE9").Value = 31665
$ws.Range("AF9").Value = 0.66
$ws.Range("AG9").Value = 575
$ws.Range("AI9").Value = 18.79
